# "Got a decent result using the catagory -> regression attempt"
#
# The author narrowed the existing AutoFilter on the "api" sheet:
#   - column A ("method") is restricted from all four statcast_pitcher*
#     variants down to just "statcast_pitcher"
#   - a new filter is added on column K ("first_year") restricting it to
#     the values 2000, 2004, 2006, 2008, 2015, 2016
# Re-applying AutoFilter recomputes which data rows are hidden.
#
# The view was also scrolled/re-selected while reviewing the filtered data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:N916")

# xlFilterValues = 7
$rng.AutoFilter(1, @("statcast_pitcher"), 7)
$rng.AutoFilter(11, @("2000", "2004", "2006", "2008", "2015", "2016"), 7)

# Re-establish the frozen header row (freeze-panes is keyed off the
# selection at the moment it's (re)enabled) and leave the selection where
# the author ended up after scrolling through the now-smaller result set.
$window = $excel.ActiveWindow
$window.FreezePanes = $false
$ws.Range("A2").Select()
$window.FreezePanes = $true

$ws.Range("B854:B855").Select()
